$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.990.51'
$ws.Range('E2').Value = '  +1.45%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.674.83'
$ws.Range('E3').Value = '  +3.20%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '537.31'
$ws.Range('E5').Value = '  +1.52%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.92'
$ws.Range('E6').Value = '  +4.32%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  +1.39%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.675.71'
$ws.Range('E9').Value = '  +2.68%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.67'
$ws.Range('E10').Value = '  +3.42%  '

$ws.Range('E11').Value = '  +2.33%  '

$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('E13').Value = '  -1.13%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.134.32'
$ws.Range('E14').Value = '  +2.49%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.915.09'
$ws.Range('E15').Value = '  +1.42%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.35'
$ws.Range('E16').Value = '  +4.13%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.665.52'
$ws.Range('E17').Value = '  +2.08%  '

$ws.Range('E18').Value = '  +1.47%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '345.45'
$ws.Range('E19').Value = '  -0.62%  '

$ws.Range('E20').Value = '  +2.24%  '

$ws.Range('E21').Value = '  +1.86%  '

$ws.Range('E22').Value = '  -0.71%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '67.59'
$ws.Range('E24').Value = '  +0.32%  '

$ws.Range('E25').Value = '  +2.76%  '

$ws.Range('E26').Value = '  -0.59%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.05%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.35'
$ws.Range('E28').Value = '  +2.72%  '

$ws.Range('E29').Value = '  +2.33%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('E31').Value = '  +3.24%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.93'
$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.15'
$ws.Range('E33').Value = '  +1.86%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '150.51'
$ws.Range('E34').Value = '  +1.04%  '

$ws.Range('E35').Value = '  +1.74%  '

$ws.Range('E36').Value = '  +3.44%  '

$ws.Range('E37').Value = '  +0.50%  '

$ws.Range('E38').Value = '  +1.76%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.827'
$ws.Range('E39').Value = '  +0.17%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '291.56'
$ws.Range('E40').Value = '  +8.20%  '

$ws.Range('E41').Value = '  +2.43%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.17%  '

$ws.Range('E43').Value = '  +1.69%  '

$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0957'
$ws.Range('E45').Value = '  -0.22%  '

$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0539'
$ws.Range('E46').Value = '  +3.86%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.979.86'
$ws.Range('E47').Value = '  +0.96%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  +2.53%  '

$ws.Range('E49').Value = '  -0.27%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.55'
$ws.Range('E50').Value = '  +1.76%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '110.00'
$ws.Range('E51').Value = '  -1.50%  '

